$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 and 6: USDC and BNB swapped positions
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "'1.011"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'309.70"
$ws.Range("E6").Value = "  +0.48%  "

# Remaining rows: price and volume(1h) updates
$ws.Range("D2").Value = "26.951.34"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "1.847.37"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D7").Value = "'0.4760"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "'0.07221"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").Value = "'0.9270"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "'0.07709"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "1.823.76"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "'6.412"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "'88.79"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "'1.014"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'0.000008636"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "26.977.86"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").Value = "'14.54"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "'5.058"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'10.66"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").Value = "'1.931"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").Value = "'152.54"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'18.19"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "'1.993"
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "'114.14"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "'4.958"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "'0.08873"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("D31").Value = "'3.324"
$ws.Range("E31").Value = "  +5.57%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").Value = "'0.7435"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "'2.716"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("D36").Value = "'1.117"
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").Value = "'0.05268"
$ws.Range("E38").Value = "  +2.65%  "
$ws.Range("D39").Value = "'2.985"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").Value = "'0.5190"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").Value = "'6.990"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").Value = "'0.1509"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'8.192"
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").Value = "'10.55"
$ws.Range("E44").Value = "  +5.65%  "
$ws.Range("D45").Value = "'0.4727"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").Value = "'1.012"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "'101.65"
$ws.Range("E47").Value = "  +3.63%  "
$ws.Range("D48").Value = "'1.601"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("D49").Value = "'65.43"
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("D50").Value = "'0.06024"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").Value = "'0.8869"
$ws.Range("E51").Value = "  +4.32%  "
